# Adding the changes we made on may 9th
# - Insert 2 new data rows right after the header (new rows 2-3), pushing the
#   previously existing 20 data rows down to rows 4-23.
# - Append 8 new data rows at the end (rows 24-31).
# Final sheet data spans A1:C31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Shift the existing 20 data rows (rows 2..21) down by 2 rows (to rows 4..23).
#    Walk from the bottom up so we never overwrite a row before it has been read.
for ($r = 21; $r -ge 2; $r--) {
    $destRow = $r + 2
    $ws.Range("A$destRow").Value2 = $ws.Range("A$r").Value2
    $ws.Range("B$destRow").Value2 = $ws.Range("B$r").Value2
    $ws.Range("C$destRow").Value2 = $ws.Range("C$r").Value2
}

# 2) Fill in the two newly-inserted rows (2 and 3).
$topRows = @(
    @(0.0167987942695617, 0.1464549452066421, 0.5288565754890442),
    @(-0.2842050492763519, -2.60824179649353, -0.3984368443489074)
)

$r = 2
foreach ($row in $topRows) {
    $ws.Range("A$r").Value2 = $row[0]
    $ws.Range("B$r").Value2 = $row[1]
    $ws.Range("C$r").Value2 = $row[2]
    $r++
}

# 3) Append the eight new rows at the bottom (rows 24..31).
$bottomRows = @(
    @(-0.2935207486152649, -1.132696866989136, 0.0025961773935705),
    @(-0.1869247704744339, -0.6145304441452026, -0.2397646158933639),
    @(-0.1950187236070633, 0.0448985956609249, -0.3084869384765625),
    @(0.0444404482841491, 0.1437060534954071, -0.3220787048339844),
    @(0.2638937830924988, 0.3306308090686798, -0.3280346393585205),
    @(-0.0166460778564214, 0.1469130963087082, 0.2872593700885772),
    @(-0.1206458881497383, 0.0920879393815994, 0.0424551330506801),
    @(0.0326812900602817, -0.0403171069920063, 0.0131336031481623)
)

$r = 24
foreach ($row in $bottomRows) {
    $ws.Range("A$r").Value2 = $row[0]
    $ws.Range("B$r").Value2 = $row[1]
    $ws.Range("C$r").Value2 = $row[2]
    $r++
}
